$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29: update title
$ws.Range("D29").Value = "[PyTorch] 머신러닝, 딥러닝 프로젝트 설계하고 템플릿 구성하기"

# Row 44: update title and link
$ws.Range("D44").Value = "Mobile Access Edge Computing (MEC) - 3GPP"
$ws.Range("E44").Value = "https://engineering-ladder.tistory.com/93"

# Row 50: update title and link
$ws.Range("D50").Value = "VQ-VAE"
$ws.Range("E50").Value = "http://incredible.egloos.com/7526008"
